$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2').Value = 'backup@backdoor.com, System, system'
$ws.Range('G3').Value = 'dnasr281@gmail.com, System'
$ws.Range('G5').Value = 'backup@backdoor.com, System'
$ws.Range('G6').Value = 'dnasr281@gmail.com, System'
$ws.Range('G8').Value = 'backup@backdoor.com, System'
$ws.Range('G10').Value = 'dnasr281@gmail.com, System'
$ws.Range('G11').Value = 'dnasr281@gmail.com, System'
$ws.Range('G12').Value = 'dnasr281@gmail.com, System'
$ws.Range('G13').Value = 'dnasr281@gmail.com, System'
$ws.Range('G14').Value = 'dnasr281@gmail.com, System'
$ws.Range('G15').Value = 'dnasr281@gmail.com, System'
$ws.Range('G17').Value = 'dnasr281@gmail.com, System'
$ws.Range('G18').Value = 'dnasr281@gmail.com, System'
$ws.Range('G19').Value = 'dnasr281@gmail.com, System'
$ws.Range('G20').Value = 'dnasr281@gmail.com, System'
$ws.Range('G21').Value = 'dnasr281@gmail.com, System'
$ws.Range('G22').Value = 'dnasr281@gmail.com, System'
$ws.Range('G26').Value = 'dnasr281@gmail.com, System'
$ws.Range('G28').Value = 'backup@backdoor.com, System, system'
$ws.Range('G29').Value = 'dnasr281@gmail.com, System'
$ws.Range('G31').Value = 'backup@backdoor.com, System'
$ws.Range('G32').Value = 'dnasr281@gmail.com, System'
$ws.Range('G34').Value = 'backup@backdoor.com, System'
$ws.Range('G36').Value = 'dnasr281@gmail.com, System'
$ws.Range('G37').Value = 'dnasr281@gmail.com, System'
$ws.Range('G38').Value = 'dnasr281@gmail.com, System'
$ws.Range('G39').Value = 'dnasr281@gmail.com, System'
$ws.Range('G40').Value = 'dnasr281@gmail.com, System'
$ws.Range('G41').Value = 'dnasr281@gmail.com, System'
$ws.Range('G43').Value = 'dnasr281@gmail.com, System'
$ws.Range('G44').Value = 'dnasr281@gmail.com, System'
$ws.Range('G45').Value = 'dnasr281@gmail.com, System'
$ws.Range('G46').Value = 'dnasr281@gmail.com, System'
$ws.Range('G47').Value = 'dnasr281@gmail.com, System'
$ws.Range('G48').Value = 'dnasr281@gmail.com, System'
$ws.Range('G50').Value = 'dnasr281@gmail.com, System'
$ws.Range('G52').Value = 'dnasr281@gmail.com, System'
$ws.Range('G54').Value = 'backup@backdoor.com, System, system'
$ws.Range('G55').Value = 'dnasr281@gmail.com, System'
$ws.Range('G57').Value = 'backup@backdoor.com, System'
$ws.Range('G58').Value = 'dnasr281@gmail.com, System'
$ws.Range('G60').Value = 'backup@backdoor.com, System'
$ws.Range('G62').Value = 'dnasr281@gmail.com, System'
$ws.Range('G63').Value = 'dnasr281@gmail.com, System'
$ws.Range('G64').Value = 'dnasr281@gmail.com, System'
$ws.Range('G65').Value = 'dnasr281@gmail.com, System'
$ws.Range('G66').Value = 'dnasr281@gmail.com, System'
$ws.Range('G67').Value = 'dnasr281@gmail.com, System'
$ws.Range('G69').Value = 'dnasr281@gmail.com, System'
$ws.Range('G70').Value = 'dnasr281@gmail.com, System'
$ws.Range('G71').Value = 'dnasr281@gmail.com, System'
$ws.Range('G72').Value = 'dnasr281@gmail.com, System'
$ws.Range('G73').Value = 'dnasr281@gmail.com, System'
$ws.Range('G74').Value = 'dnasr281@gmail.com, System'
$ws.Range('G76').Value = 'dnasr281@gmail.com, System'
$ws.Range('G78').Value = 'dnasr281@gmail.com, System'
$ws.Range('G80').Value = 'backup@backdoor.com, System'
$ws.Range('G81').Value = 'backup@backdoor.com, System'
$ws.Range('G82').Value = 'backup@backdoor.com, System'
$ws.Range('G83').Value = 'dnasr281@gmail.com, System'
$ws.Range('G84').Value = 'dnasr281@gmail.com, System'
$ws.Range('G85').Value = 'dnasr281@gmail.com, System'
$ws.Range('G86').Value = 'dnasr281@gmail.com, System'
$ws.Range('G87').Value = 'dnasr281@gmail.com, admin@admin.com'
$ws.Range('G90').Value = 'dnasr281@gmail.com, System'
$ws.Range('G92').Value = 'dnasr281@gmail.com, System'
$ws.Range('G93').Value = 'dnasr281@gmail.com, System'
$ws.Range('G94').Value = 'dnasr281@gmail.com, System'
$ws.Range('G96').Value = 'dnasr281@gmail.com, System'
$ws.Range('G99').Value = 'dnasr281@gmail.com, System'
$ws.Range('G101').Value = 'dnasr281@gmail.com, System'
$ws.Range('G106').Value = 'backup@backdoor.com, System'
$ws.Range('G107').Value = 'backup@backdoor.com, System'
$ws.Range('G108').Value = 'backup@backdoor.com, System'
$ws.Range('G109').Value = 'dnasr281@gmail.com, System'
$ws.Range('G110').Value = 'dnasr281@gmail.com, System'
$ws.Range('G111').Value = 'dnasr281@gmail.com, System'
$ws.Range('G112').Value = 'dnasr281@gmail.com, System'
$ws.Range('G113').Value = 'dnasr281@gmail.com, admin@admin.com'
$ws.Range('G116').Value = 'dnasr281@gmail.com, System'
$ws.Range('G118').Value = 'dnasr281@gmail.com, System'
$ws.Range('G119').Value = 'dnasr281@gmail.com, System'
$ws.Range('G120').Value = 'dnasr281@gmail.com, System'
$ws.Range('G122').Value = 'dnasr281@gmail.com, System'
$ws.Range('G125').Value = 'dnasr281@gmail.com, System'
$ws.Range('G127').Value = 'dnasr281@gmail.com, System'
$ws.Range('G132').Value = 'backup@backdoor.com, System'
$ws.Range('G133').Value = 'backup@backdoor.com, System'
$ws.Range('G134').Value = 'backup@backdoor.com, System'
$ws.Range('G135').Value = 'dnasr281@gmail.com, System'
$ws.Range('G136').Value = 'dnasr281@gmail.com, System'
$ws.Range('G137').Value = 'dnasr281@gmail.com, System'
$ws.Range('G138').Value = 'dnasr281@gmail.com, System'
$ws.Range('G139').Value = 'dnasr281@gmail.com, admin@admin.com'
$ws.Range('G142').Value = 'dnasr281@gmail.com, System'
$ws.Range('G144').Value = 'dnasr281@gmail.com, System'
$ws.Range('G145').Value = 'dnasr281@gmail.com, System'
$ws.Range('G146').Value = 'dnasr281@gmail.com, System'
$ws.Range('G148').Value = 'dnasr281@gmail.com, System'
$ws.Range('G151').Value = 'dnasr281@gmail.com, System'
$ws.Range('G153').Value = 'dnasr281@gmail.com, System'
